$p = $ppt.ActivePresentation
$s = $p.Slides.Item(30)
$tf = $s.Shapes.Item(2).TextFrame
$tr = $tf.TextRange

# Paragraph (lvl 1): "if the receiver did not receive a segment " -> "if the receiver did not receive a packet"
$para3 = $tr.Paragraphs(3)
$run3 = $para3.Runs(1)
$run3.Text = "if the receiver did not receive a packet"

# Paragraph (lvl 1): "but did receive a subsequent few segments (..." ->
#   split into "but did receive a subsequent " + "few packets ("  (rest of paragraph unchanged)
$para4 = $tr.Paragraphs(4)
$splitPart = $para4.Characters(30, 14)
$splitPart.Text = "few packets ("
